# Update the "Förändrad" (changed) date column (C) for rows 2-17
# from serial date 45185 (2023-09-16) to 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
